# Set column K (Ανταπόκριση / Response) from 'yes' to 'no' for the specific
# rows identified in the commit diff (knn model / best-k overfitting fix).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(16, 18, 21, 23, 39, 48, 49, 62, 89, 90, 95, 100, 101, 102, 104, 106, 112, 119, 121, 129, 135, 144, 166, 173, 174, 177, 179, 182, 184, 202, 214, 219, 225, 229, 234, 270, 272, 275, 284, 290, 295, 308, 311, 315, 323, 329, 332, 336, 338, 351, 363, 371, 377, 381, 384, 396, 405, 411, 413, 419, 430, 439, 448, 451, 458, 468, 479, 509, 512, 513, 518, 521, 524, 527, 529, 530, 539, 549, 550, 551, 553, 564, 566, 568, 569, 570, 571, 577, 578, 589, 591, 626, 628, 648, 652, 658, 659, 678, 679, 696, 697, 698, 704, 711, 726, 727, 728, 739, 740, 778, 794, 810, 816, 823, 826, 833, 839, 840, 850, 853, 857, 877, 882, 893, 900, 919, 920, 964, 975, 981, 983, 990)

foreach ($r in $rows) {
    $ws.Range("K$r").Value = "no"
}

